$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1070.53
$ws.Range("C3").Value = 1072.12
$ws.Range("C4").Value = 1012.18
$ws.Range("C5").Value = 1019.23
$ws.Range("C6").Value = 1019.23
